$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Info")

$ws.Range("C8").Value = -1.98
$ws.Range("C9").Value = -1.33
$ws.Range("C12").Formula = "=100-56.3"
$ws.Range("C13").Value = 1.1499999999999999

$ws.Activate()
$ws.Range("C12:C13").Select()
